# Actualización automática 2025-07-17 15:55:08
# Target sheet: "CUMPLIMIENTO MENSUAL" (3rd worksheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Update VENTA (D) figures for row 3 (PORCELANATO) and row 4 (TOTAL);
# POR CUMPLIR (E) and CUMPLIMIENTO (F) are recomputed accordingly.
$ws.Range("D3").Value = 5114.65
$ws.Range("E3").Value = 8608.690000000001
$ws.Range("F3").Value = 0.3726971713883063

$ws.Range("D4").Value = 10161.17
$ws.Range("E4").Value = 3562.17
$ws.Range("F4").Value = 0.7404298079039068

# Narrow column F slightly (width 25 -> 24 in OOXML character units)
$ws.Columns.Item(6).ColumnWidth = 23.17
